# Daily auto push: 2026-02-22 03:15 UTC
# Insert a new data row for 2026/02/22 (日) at row 832, pushing the
# existing rows (the 2026/12/29 .. 2027/01/05 block) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 832 (shifts 832..873 down to 833..874
# and expands the sheet dimension to A1:D874 automatically).
$ws.Rows(832).Insert()

# Force column A to plain text first so the date-like string "2026/02/22" is
# stored as literal text (matching the rest of the sheet) instead of being
# auto-converted into a date serial number.
$ws.Cells.Item(832, 1).NumberFormat = "@"
$ws.Cells.Item(832, 1).Value() = "2026/02/22"
# Strip any formatting picked up from the NumberFormat change above so the new
# row stays unstyled, just like every other data row in the sheet.
$ws.Cells.Item(832, 1).Style = "Normal"

$ws.Cells.Item(832, 2).Value() = "日"
$ws.Cells.Item(832, 3).Value() = 8
$ws.Cells.Item(832, 4).Value() = 42
